# Replace embedded line breaks in several BrandName / Manufacturer cells
# with plain spaces across the workbook's sheets.
$wb = $excel.ActiveWorkbook

# "Adult Vaccine " sheet
$ws = $wb.Worksheets.Item("Adult Vaccine ")
$ws.Range("B12").Value = "Tetanus  Diphtheria Toxoids Adsorbed for Adults No Preservative"

# "Pediatric Influenza Vaccine " sheet
$ws = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$ws.Range("B3").Value  = "Fluzone Pediatric dose No Preservative"
$ws.Range("B6").Value  = "Fluarix Preservative-Free"
$ws.Range("B9").Value  = "FluMist No Preservative"
$ws.Range("B10").Value = "Afluria No Preservative"
$ws.Range("H10").Value = "Merck (CSL product)"
$ws.Range("H11").Value = "Merck (CSL product)"
$ws.Range("B12").Value = "Afluria No Preservative"
$ws.Range("H12").Value = "Merck (CSL product)"

# "Adult Influenza Vaccine " sheet
$ws = $wb.Worksheets.Item("Adult Influenza Vaccine ")
$ws.Range("B5").Value  = "Agriflu No Preservative"
$ws.Range("B7").Value  = "Fluvirin Preservative-free"
$ws.Range("B8").Value  = "Fluraix Preservative-free"
$ws.Range("B10").Value = "Flumist No Preservative"
